$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("lambda")
$ws1.Range("G3").Value = 0.1204325429604989
$ws1.Range("J3").Value = 0.2874501708223065
$ws1.Range("T3").Value = 0.5921172862171945
$ws1.Range("K5").Value = 0
$ws1.Range("I6").Value = 0
$ws1.Range("K6").Value = 0
$ws1.Range("X6").Value = 0
$ws1.Range("C7").Value = 0
$ws1.Range("H7").Value = 0
$ws1.Range("I7").Value = 0
$ws1.Range("K7").Value = 0
$ws1.Range("M7").Value = 0
$ws1.Range("Q7").Value = 0
$ws1.Range("S7").Value = 0
$ws1.Range("X7").Value = 0
$ws1.Range("G8").Value = 0.09475633363640217
$ws1.Range("J8").Value = 0.8675622404910103
$ws1.Range("O8").Value = 0.03768142587258756
$ws1.Range("F9").Value = 0.9492480771877249
$ws1.Range("G9").Value = 0.05075192281227514
$ws1.Range("C10").Value = 0
$ws1.Range("H10").Value = 0
$ws1.Range("M10").Value = 0
$ws1.Range("Q10").Value = 0
$ws1.Range("S10").Value = 0
$ws1.Range("X10").Value = 0
$ws1.Range("Y10").Value = 0
$ws1.Range("E11").Value = 0.4406500996151126
$ws1.Range("F11").Value = 0.5101026577201669
$ws1.Range("G11").Value = 0.03717187094160868
$ws1.Range("T11").Value = 0.01207537172311179
$ws1.Range("O12").Value = 0.08433330658934152
$ws1.Range("T12").Value = 0.5367373836391109
$ws1.Range("U12").Value = 0.3789293097715476
$ws1.Range("G13").Value = 0.1003814589172418
$ws1.Range("J13").Value = 0.2714521275478711
$ws1.Range("T13").Value = 0.6281664135348871
$ws1.Range("H15").Value = 0
$ws1.Range("L15").Value = 0
$ws1.Range("R15").Value = 0
$ws1.Range("Y15").Value = 0
$ws1.Range("G17").Value = 0.1123218570984662
$ws1.Range("J17").Value = 0.2543759704877028
$ws1.Range("T17").Value = 0.633302172413831
$ws1.Range("O18").Value = 0.2437926068535217
$ws1.Range("T18").Value = 0.4756442191313157
$ws1.Range("U18").Value = 0.2805631740151626
$ws1.Range("G19").Value = 0.06163052902152446
$ws1.Range("J19").Value = 0.1395750216075701
$ws1.Range("T19").Value = 0.7987944493709055
$ws1.Range("C20").Value = 0
$ws1.Range("K20").Value = 0
$ws1.Range("L20").Value = 0
$ws1.Range("M20").Value = 0
$ws1.Range("Q20").Value = 0
$ws1.Range("R20").Value = 0
$ws1.Range("S20").Value = 0
$ws1.Range("X20").Value = 0
$ws1.Range("L21").Value = 0
$ws1.Range("R21").Value = 0
$ws1.Range("F24").Value = 0.1287728355876866
$ws1.Range("G24").Value = 0.05347445417152326
$ws1.Range("J24").Value = 0.3990366189573973
$ws1.Range("T24").Value = 0.4187160912833928
$ws1.Range("J25").Value = 0.8870902668099052
$ws1.Range("O25").Value = 0.1129097331900948

$ws2 = $wb.Worksheets.Item("omega")
$ws2.Range("B2").Value = 0
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 0
$ws2.Range("H2").Value = 0
$ws2.Range("I2").Value = 0
$ws2.Range("J2").Value = 0
$ws2.Range("K2").Value = 0
$ws2.Range("L2").Value = 0
$ws2.Range("M2").Value = 0
$ws2.Range("N2").Value = 0
$ws2.Range("O2").Value = 0.4015376883739908
$ws2.Range("P2").Value = 0
$ws2.Range("Q2").Value = 0
$ws2.Range("R2").Value = 0
$ws2.Range("S2").Value = 0.5984623116260092
$ws2.Range("T2").Value = 0
$ws2.Range("U2").Value = 0
$ws2.Range("V2").Value = 0
$ws2.Range("W2").Value = 0
$ws2.Range("X2").Value = 0
$ws2.Range("Y2").Value = 0
$ws2.Range("B3").Value = 0.6599006090483395
$ws2.Range("E3").Value = 0.3400993909516606
$ws2.Range("O4").Value = 0.5808446626520776
$ws2.Range("S4").Value = 0.4191553373479224
$ws2.Range("C5").Value = 0
$ws2.Range("F5").Value = 0
$ws2.Range("G5").Value = 0
$ws2.Range("H5").Value = 0
$ws2.Range("I5").Value = 0
$ws2.Range("K5").Value = 0
$ws2.Range("L5").Value = 0
$ws2.Range("M5").Value = 0
$ws2.Range("N5").Value = 0
$ws2.Range("P5").Value = 0
$ws2.Range("Q5").Value = 0
$ws2.Range("R5").Value = 0
$ws2.Range("S5").Value = 0
$ws2.Range("T5").Value = 0
$ws2.Range("V5").Value = 0
$ws2.Range("W5").Value = 0
$ws2.Range("X5").Value = 0
$ws2.Range("B6").Value = 0.03009285297811815
$ws2.Range("E6").Value = 0.9699071470218819
$ws2.Range("B7").Value = 0.2870887548328691
$ws2.Range("E7").Value = 0.7129112451671309
$ws2.Range("B8").Value = 0.973185466738121
$ws2.Range("E8").Value = 0.02681453326187897
$ws2.Range("B9").Value = 0.04313588915212152
$ws2.Range("E9").Value = 0.9568641108478785
$ws2.Range("O10").Value = 0.4236104703796881
$ws2.Range("S10").Value = 0.5763895296203119
$ws2.Range("B11").Value = 0.03278614440816596
$ws2.Range("E11").Value = 0.967213855591834
$ws2.Range("B12").Value = 0.9012589160324767
$ws2.Range("E12").Value = 0.09874108396752337
$ws2.Range("B13").Value = 0.6579942322400196
$ws2.Range("E13").Value = 0.3420057677599804
$ws2.Range("B14").Value = 0.8150927473684779
$ws2.Range("E14").Value = 0.1849072526315221
$ws2.Range("D15").Value = 0
$ws2.Range("J15").Value = 0
$ws2.Range("U15").Value = 0
$ws2.Range("Y15").Value = 0
$ws2.Range("B16").Value = 0.5778522339386156
$ws2.Range("E16").Value = 0.4221477660613844
$ws2.Range("B17").Value = 0.6468546256285614
$ws2.Range("E17").Value = 0.3531453743714386
$ws2.Range("B18").Value = 0.9985062376000774
$ws2.Range("E18").Value = 0.001493762399922614
$ws2.Range("B19").Value = 0.6077264995218093
$ws2.Range("E19").Value = 0.3922735004781907
$ws2.Range("B20").Value = 0.5601545142717921
$ws2.Range("E20").Value = 0.4398454857282079
$ws2.Range("B21").Value = 0.2861161368141096
$ws2.Range("O21").Value = 0.5244625113925655
$ws2.Range("S21").Value = 0.1894213517933249
$ws2.Range("B22").Value = 0.8334069326097517
$ws2.Range("E22").Value = 0.1665930673902483
$ws2.Range("B23").Value = 0.1439498693142103
$ws2.Range("E23").Value = 0.8560501306857897
$ws2.Range("B24").Value = 0.6614140197534236
$ws2.Range("E24").Value = 0.3385859802465764
$ws2.Range("O25").Value = 0.4886904583826819
$ws2.Range("S25").Value = 0.5113095416173181
